{"js": "// Word Rules section grammar/content fixes:\n//   1) \"The game is a two-player game set on a 8x8 ...\" -> \"... on an 8x8 ...\"\n//      (fixes the \"a\" -> \"an\" article before the vowel sound in \"8x8\")\n//   2) \"Each player has 12 game pieces: 4 guards, 4 mercenaries, ...\"\n//      -> \"... 4 guards, 5 mercenaries, ...\"\n\nconst body = context.document.body;\n\n// --- 1) \"a\" -> \"an\" before \"8x8\" -------------------------------------------\n// Scope the search to the paragraph that contains the sentence so we hit the\n// correct \"a\" (the word \"a\" also appears earlier in \"is a two-player\").\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet rulesParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"two-player game set on a 8x8\") !== -1) {\n    rulesParagraph = p;\n    break;\n  }\n}\n\nif (rulesParagraph) {\n  const articleMatches = rulesParagraph.search(\"a\", {\n    matchCase: true,\n    matchWholeWord: true\n  });\n  articleMatches.load(\"text\");\n  await context.sync();\n\n  // The last whole-word \"a\" in this sentence is the article right before\n  // \"8x8\" (\"The game is [a] two-player game set on [a] 8x8 ...\").\n  if (articleMatches.items.length > 0) {\n    const article = articleMatches.items[articleMatches.items.length - 1];\n    article.insertText(\"an\", \"Replace\");\n    await context.sync();\n  }\n}\n\n// --- 2) \"4 mercenaries\" -> \"5 mercenaries\" ----------------------------------\nconst mercenaryMatches = body.search(\"4 mercenaries\", { matchCase: true });\nmercenaryMatches.load(\"text\");\nawait context.sync();\n\nif (mercenaryMatches.items.length > 0) {\n  mercenaryMatches.items[0].insertText(\"5 mercenaries\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word Rules section grammar/content fixes:\n#   1) \"The game is a two-player game set on a 8x8 ...\" -> \"... on an 8x8 ...\"\n#      (fixes the \"a\" -> \"an\" article before the vowel sound in \"8x8\")\n#   2) \"Each player has 12 game pieces: 4 guards, 4 mercenaries, ...\"\n#      -> \"... 4 guards, 5 mercenaries, ...\"\n\n$d = $word.ActiveDocument\n\n# --- 1) \"a\" -> \"an\" before \"8x8\" --------------------------------------------\n# Locate the end of \"set on \" (unique in the document) and grab the single\n# character right after it -- that is the article \"a\" that needs to become\n# \"an\". This avoids touching the separately-styled \"8\", \"x\", \"8\" runs that\n# follow it.\n$locator = $d.Content\n$locator.Find.ClearFormatting()\n$locator.Find.Text = \"set on \"\n$locator.Find.Forward = $true\n$locator.Find.Wrap = 0\n$foundLocator = $locator.Find.Execute()\n\nif ($foundLocator) {\n    $articleStart = $locator.End\n    $articleRange = $d.Range($articleStart, $articleStart + 1)\n    if ($articleRange.Text -eq \"a\") {\n        $articleRange.Text = \"an\"\n    }\n}\n\n# --- 2) \"4 mercenaries\" -> \"5 mercenaries\" ----------------------------------\n$mercRange = $d.Content\n$mercRange.Find.ClearFormatting()\n$mercRange.Find.Execute(\"4 mercenaries\", $false, $false, $false, $false, $false, $true, 1, $false, \"5 mercenaries\", 2)\n"}
